$d = $word.ActiveDocument

$target = "El Administrador da clic al bot" + [char]0x00F3 + "n Crear horario."

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $r = $p.Range
    $t = $r.Text.TrimEnd([char]13, [char]7).Trim()
    if ($t -eq $target) {
        # Delete the whole paragraph, including its paragraph mark,
        # so the following paragraph collapses up without leaving a blank line.
        $r.Delete()
        break
    }
}
